$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "invest"
$ws.Range("C2").Value = -0.3228

# Row 4
$ws.Range("B4").Value = "uncertain"
$ws.Range("C4").Value = -0.3122

# Row 7
$ws.Range("B7").Value = "inflation"
$ws.Range("C7").Value = -0.0236

# Row 8
$ws.Range("B8").Value = "trade"
$ws.Range("C8").Value = 0.6127

# Row 9
$ws.Range("B9").Value = "interest"
$ws.Range("C9").Value = 0.07969999999999999

# Row 10
$ws.Range("B10").Value = "uncertain"
$ws.Range("C10").Value = -0.4164

# Row 11
$ws.Range("B11").Value = "invest"
$ws.Range("C11").Value = 0.4366

# Row 12
$ws.Range("B12").Value = "trade"
$ws.Range("C12").Value = -0.127

# Row 13
$ws.Range("B13").Value = "uncertain"
$ws.Range("C13").Value = 0.06569999999999999

# Row 14
$ws.Range("B14").Value = "interest"
$ws.Range("C14").Value = -0.0998

# Row 15
$ws.Range("B15").Value = "invest"
$ws.Range("C15").Value = -0.1474

# Row 17
$ws.Range("B17").Value = "trade"
$ws.Range("C17").Value = -0.1159

# Row 18
$ws.Range("B18").Value = "inflation"
$ws.Range("C18").Value = 0.0035

# Row 19
$ws.Range("B19").Value = "interest"
$ws.Range("C19").Value = -0.0366

# Row 20
$ws.Range("B20").Value = "invest"
$ws.Range("C20").Value = 0.1352

# Row 21
$ws.Range("B21").Value = "uncertain"
$ws.Range("C21").Value = -0.2433
